$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Statistical " + "Analysis:" -> single run "Statistical Analysis:"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Statistical Analysis:", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Statistical Analysis:", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) " - Research" + " Intern" -> single run " - Research Intern"
#    (search only the text of the two runs being merged so the
#    preceding "Oak Ridge National Laboratory" run is left untouched)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" Intern", $false, $false, $false, `
    $false, $false, $true, 1, $false, " Intern", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Shrink the run of 14 spaces right before "April" down to 13
#    spaces, without letting it merge with the identically formatted
#    "April" run. Any text edit to that run would normally coalesce it
#    with format-identical neighbours, so the "April" run's formatting
#    is perturbed (via a property that round-trips byte-for-byte) just
#    long enough to make the two runs differ, then restored.
# ---------------------------------------------------------------------
$april = $d.Content
$april.Find.Execute("April") | Out-Null
$aprilStart = $april.Start
$aprilEnd = $april.End
$origColor = $april.Font.Color
$april.Font.Color = 255

$spaces = $d.Range($aprilStart - 14, $aprilStart)
$spaces.Text = "             "

# "April"'s Range object is now stale (one char removed ahead of it),
# so re-derive its new position arithmetically instead of trusting it.
$april2 = $d.Range($aprilStart - 1, $aprilEnd - 1)
$april2.Font.Color = $origColor

# ---------------------------------------------------------------------
# 4) "Parametric nasophary" + "n" + "geal swab for sampling COVID-19
#    and other respiratory viruses" -> single run
# ---------------------------------------------------------------------
$oldCovid = "Parametric nasophary" + "n" + "geal swab for sampling COVID-19 and other respiratory viruses"
$newCovid = "Parametric nasopharyngeal swab for sampling COVID-19 and other respiratory viruses"
$d.Content.Find.Execute($oldCovid, $false, $false, $false, `
    $false, $false, $true, 1, $false, $newCovid, 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Mark the built-in "Default Paragraph Font" style as semi-hidden.
# ---------------------------------------------------------------------
try {
    $d.Styles("Default Paragraph Font").Hidden = $true
} catch {
}
